# Updates the crypto price/volume table to reflect the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.025.21"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "2.447.35"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.32"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.95"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.29"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "2.830.62"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "2.451.27"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "45.908.32"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.33"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.06"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.99"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.01"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +6.49%  "
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0761"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.97"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.02"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "1.959.28"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("E48").Value = "  +8.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.03"
$ws.Range("E49").Value = "  -7.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.64"
$ws.Range("E50").Value = "  +5.08%  "
$ws.Range("E51").Value = "  +5.63%  "
